# Auto-generated script to apply cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.694.25"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +5.97%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.735.10"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +4.82%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.68"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +3.90%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5456"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +3.74%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.08%  "

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.17%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06720"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +5.37%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.89"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +6.37%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07789"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.28%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.699"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.82%  "

# Row 13
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.973.68"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +4.80%  "

# Row 14
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.715.98"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.40%  "

# Row 15
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +6.22%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8424"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.97%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.22"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +5.35%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "27.702.47"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +6.07%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "228.21"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +19.84%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.833"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.93%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.004"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.03%  "

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +5.16%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.228"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +3.73%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.004"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.02%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.24"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.51%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.723"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +13.24%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1251"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +4.03%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.477"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.80%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "17.13"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +6.95%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05706"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.92%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.314"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.83%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.702"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +5.87%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.521"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.95%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.689"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +6.71%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9764"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.87%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.854"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.95%  "

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.27%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5992"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.61%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01671"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +4.51%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.939"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.61%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8514"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.81%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.049.46"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.54%  "

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.01%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.68"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.26%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.878.26"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +4.70%  "

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +10.01%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "59.68"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.328"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.34%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4433"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.12%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05335"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.14%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9986"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.54%  "

Write-Host "Applied 97 cell updates"